# Apply cryptos list update (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.566.34'
$ws.Cells.Item(2, 5).Value = '  +0.70%  '
$ws.Cells.Item(3, 4).Value = '3.447.63'
$ws.Cells.Item(3, 5).Value = '  +1.94%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '581.66'
$ws.Cells.Item(5, 5).Value = '  +1.28%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '146.26'
$ws.Cells.Item(6, 5).Value = '  +6.38%  '
$ws.Cells.Item(7, 4).Value = '3.447.55'
$ws.Cells.Item(7, 5).Value = '  +1.94%  '
$ws.Cells.Item(8, 5).Value = '  +0.06%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.475'
$ws.Cells.Item(9, 5).Value = '  +1.28%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '7.63'
$ws.Cells.Item(10, 5).Value = '  -0.17%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.126'
$ws.Cells.Item(11, 5).Value = '  +2.83%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.390'
$ws.Cells.Item(12, 5).Value = '  +2.38%  '
$ws.Cells.Item(13, 4).Value = '4.040.68'
$ws.Cells.Item(13, 5).Value = '  +2.00%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '27.93'
$ws.Cells.Item(14, 5).Value = '  +8.27%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0000174'
$ws.Cells.Item(16, 5).Value = '  +1.03%  '
$ws.Cells.Item(17, 4).Value = '3.426.60'
$ws.Cells.Item(17, 5).Value = '  +1.18%  '
$ws.Cells.Item(18, 4).Value = '61.754.72'
$ws.Cells.Item(18, 5).Value = '  +0.80%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.22'
$ws.Cells.Item(19, 5).Value = '  +8.30%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '14.34'
$ws.Cells.Item(20, 5).Value = '  +3.72%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '9.53'
$ws.Cells.Item(21, 5).Value = '  +2.10%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '387.95'
$ws.Cells.Item(22, 5).Value = '  +2.97%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.566'
$ws.Cells.Item(23, 5).Value = '  +2.72%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '73.68'
$ws.Cells.Item(24, 5).Value = '  +3.77%  '
$ws.Cells.Item(25, 5).Value = '  +0.08%  '
$ws.Cells.Item(26, 5).Value = '  -0.01%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.0000123'
$ws.Cells.Item(27, 5).Value = '  -1.94%  '
$ws.Cells.Item(28, 4).Value = '3.602.23'
$ws.Cells.Item(28, 5).Value = '  +2.38%  '
$ws.Cells.Item(29, 5).Value = '  +1.42%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '7.66'
$ws.Cells.Item(30, 5).Value = '  +3.27%  '
$ws.Cells.Item(31, 5).Value = '  +0.09%  '
$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '8.19'
$ws.Cells.Item(32, 5).Value = '  +1.71%  '
$ws.Cells.Item(33, 2).Value = 'Fetch.AI'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.47'
$ws.Cells.Item(33, 5).Value = '  -12.11%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.19'
$ws.Cells.Item(34, 5).Value = '  +2.09%  '
$ws.Cells.Item(35, 5).Value = '  +0.06%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '24.11'
$ws.Cells.Item(36, 5).Value = '  +3.03%  '
$ws.Cells.Item(37, 4).Value = '3.482.56'
$ws.Cells.Item(37, 5).Value = '  +2.26%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '6.99'
$ws.Cells.Item(38, 5).Value = '  +2.40%  '
$ws.Cells.Item(39, 2).Value = 'NEARProtocol'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.14'
$ws.Cells.Item(39, 5).Value = '  -0.13%  '
$ws.Cells.Item(40, 2).Value = 'ImmutableX'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.56'
$ws.Cells.Item(40, 5).Value = '  +0.72%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '167.07'
$ws.Cells.Item(41, 5).Value = '  +1.46%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0783'
$ws.Cells.Item(42, 5).Value = '  +3.17%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '27.46'
$ws.Cells.Item(43, 5).Value = '  +7.35%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.804'
$ws.Cells.Item(44, 5).Value = '  +3.68%  '
$ws.Cells.Item(45, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.00'
$ws.Cells.Item(45, 5).Value = '  +0.14%  '
$ws.Cells.Item(46, 2).Value = 'Filecoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '4.51'
$ws.Cells.Item(46, 5).Value = '  +3.75%  '
$ws.Cells.Item(47, 2).Value = 'OKB'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '42.46'
$ws.Cells.Item(47, 5).Value = '  +1.75%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.72'
$ws.Cells.Item(48, 5).Value = '  +1.07%  '
$ws.Cells.Item(49, 2).Value = 'Maker'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(49, 4).Value = '2.574.38'
$ws.Cells.Item(49, 5).Value = '  +1.10%  '
$ws.Cells.Item(50, 2).Value = 'ONDO'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.16'
$ws.Cells.Item(50, 5).Value = '  -2.55%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '6.93'
$ws.Cells.Item(51, 5).Value = '  +2.25%  '
